$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "objetos" (objects/devices) table gained a "type" column and a couple
# of extra rows, and the Alexa row grew two trailing columns (owner + a
# second boolean flag). Rewrite the whole 4-row table with the corrected
# layout: Name | Type | Value | Active  (Alexa row also carries Owner | Active2)
# ---------------------------------------------------------------------------

# Row 1: Lampada do quarto / lampada / 100 / FALSE
$ws.Range("A1").Value = "Lampada do quarto"
$ws.Range("B1").Value = "lampada"
$ws.Range("C1").Value = 100
$ws.Range("D1").Value = $false

# Row 2: Ar da sala / climatizadores / 18 / TRUE
$ws.Range("A2").Value = "Ar da sala"
$ws.Range("B2").Value = "climatizadores"
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = $true

# Row 3: Alexa / alto falante / TRUE / 50 / Leo Santana / TRUE
$ws.Range("A3").Value = "Alexa"
$ws.Range("B3").Value = "alto falante"
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = "Leo Santana"
$ws.Range("F3").Value = $true

# Row 4 (new): Lampada da sala / lampada / 0 / FALSE
$ws.Range("A4").Value = "Lampada da sala"
$ws.Range("B4").Value = "lampada"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = $false

# Remove the now-unused trailing cells left over from the old layout (column E
# used to hold data on rows 1-2; it's blank now that the columns shifted).
$ws.Range("E1").Clear()
$ws.Range("E2").Clear()

# ---------------------------------------------------------------------------
# Formatting: only the first row keeps explicit number formats / alignment
# (Name general, Type right-aligned numeric-style, Value right-aligned
# numeric-style, Active centered); the rest of the data rows use the sheet's
# plain default formatting.
# ---------------------------------------------------------------------------
$ws.Range("A1").HorizontalAlignment = 1        # xlGeneral
$ws.Range("B1").NumberFormat = "#,##0"
$ws.Range("B1").HorizontalAlignment = -4152    # xlRight
$ws.Range("C1").NumberFormat = "#,##0"
$ws.Range("C1").HorizontalAlignment = -4152    # xlRight
$ws.Range("D1").HorizontalAlignment = -4108    # xlCenter

$ws.Range("A2:D2").Style = "Normal"
$ws.Range("A3:F3").Style = "Normal"
$ws.Range("A4:D4").Style = "Normal"
